$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 12
$ws.Range("F2").Value = 4.5
$ws.Range("G2").Value = 6.8
$ws.Range("H2").Value = 1.74
$ws.Range("I2").Value = 1.93
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 1.31
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.25
$ws.Range("O2").Value = 1.31
$ws.Range("P2").Value = 1.89
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.35
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 1.97
$ws.Range("V2").Value = 2.08
$ws.Range("W2").Value = 1.24
$ws.Range("Y2").Value = 9.6
$ws.Range("Z2").Value = 12
$ws.Range("F3").Value = 1.47
$ws.Range("I3").Value = 9
$ws.Range("V3").Value = 1.13
$ws.Range("T4").Value = 1.69
$ws.Range("AD7").Value = 15.5
$ws.Range("AI7").Value = 60
$ws.Range("F7").Value = 2.56
$ws.Range("G7").Value = 2.94
$ws.Range("H7").Value = 2.68
$ws.Range("I7").Value = 3.1
$ws.Range("K7").Value = 3.75
$ws.Range("P7").Value = 1.72
$ws.Range("Q7").Value = 2.12
$ws.Range("S7").Value = 3.95
$ws.Range("U7").Value = 1.98
$ws.Range("V7").Value = 1.48
$ws.Range("W7").Value = 1.51
$ws.Range("Y7").Value = 12
$ws.Range("Z7").Value = 22
$ws.Range("W8").Value = 3.9
$ws.Range("AL10").Value = 75
$ws.Range("AN10").Value = 85
$ws.Range("F10").Value = 5.8
$ws.Range("G10").Value = 6
$ws.Range("Q10").Value = 1.85
$ws.Range("AE11").Value = 65
$ws.Range("AN11").Value = 6.2
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 1.18
$ws.Range("R11").Value = 1.67
$ws.Range("T11").Value = 1.67
$ws.Range("I12").Value = 2.22
$ws.Range("G13").Value = 1.6
$ws.Range("W13").Value = 2.68
$ws.Range("H15").Value = 2.8
$ws.Range("T15").Value = 1.68
$ws.Range("F16").Value = 1.57
$ws.Range("G16").Value = 1.64
$ws.Range("J16").Value = 4
$ws.Range("Q19").Value = 1.48
$ws.Range("W19").Value = 2.9
$ws.Range("G20").Value = 1.55
$ws.Range("P20").Value = 2.5
$ws.Range("U20").Value = 2.06
$ws.Range("G21").Value = 1.39
$ws.Range("J21").Value = 5.2
$ws.Range("N21").Value = 2.28
$ws.Range("R21").Value = 1.51
$ws.Range("W21").Value = 3.9
$ws.Range("L22").Value = 1.33
$ws.Range("S22").Value = 3.25
$ws.Range("G23").Value = 1.67
$ws.Range("AO24").Value = 5.8
$ws.Range("N24").Value = 5.4
$ws.Range("Y24").Value = 16
$ws.Range("Q25").Value = 2.1
$ws.Range("T26").Value = 1.71
$ws.Range("AE27").Value = 980
$ws.Range("AH27").Value = 15.5
$ws.Range("H27").Value = 2.88
$ws.Range("AA30").Value = 22
$ws.Range("AB30").Value = 23
$ws.Range("AE30").Value = 18
$ws.Range("AH30").Value = 17
$ws.Range("AK30").Value = 50
$ws.Range("AM30").Value = 1000
$ws.Range("AN30").Value = 1000
$ws.Range("AO30").Value = 8.6
$ws.Range("L30").Value = 1.28
$ws.Range("N30").Value = 5.3
$ws.Range("P30").Value = 2.46
$ws.Range("Q30").Value = 1.54
$ws.Range("R30").Value = 1.59
$ws.Range("T30").Value = 1.58
$ws.Range("U30").Value = 2.38
$ws.Range("AA31").Value = 260
$ws.Range("AD31").Value = 26
$ws.Range("AE31").Value = 110
$ws.Range("AH31").Value = 25
$ws.Range("AI31").Value = 100
$ws.Range("H31").Value = 7.4
$ws.Range("I31").Value = 7.6
$ws.Range("J31").Value = 4.5
$ws.Range("K31").Value = 4.6
$ws.Range("R31").Value = 1.41
$ws.Range("S31").Value = 3.35
$ws.Range("V31").Value = 1.15
$ws.Range("Z31").Value = 55
$ws.Range("O32").Value = 1.21
$ws.Range("S32").Value = 2.56
